$d = $word.ActiveDocument

# --- Change 1: merge "Add option to c" + "hange resolution by default"
# into a single run, dropping the stray _GoBack bookmark that used to
# sit at the split point. A same-text Find/Replace over the whole
# sentence causes Word to re-run the range as one run.
$d.Content.Find.Execute(
    "Add option to change resolution by default", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Add option to change resolution by default", 2) | Out-Null

# --- Change 2: fix the mis-numbered list item "5.Open configuration
# with ctrl+O." -> "6.Open configuration with ctrl+O.", which also
# relocates the _GoBack bookmark to sit right after the new "6".
$d.Content.Find.Execute(
    "5.Open configuration with ctrl+O.", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "6.Open configuration with ctrl+O.", 2) | Out-Null

# Find the paragraph that now starts with "6.Open configuration" and
# drop a point-bookmark named _GoBack right after the leading "6".
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $par = $d.Paragraphs($i)
    $t = $par.Range.Text
    if ($t.StartsWith("6.Open configuration with ctrl+O.")) {
        $pos = $par.Range.Start + 1
        $pt = $d.Range($pos, $pos)
        $d.Bookmarks.Add("_GoBack", $pt) | Out-Null
        break
    }
}
